# "função plano alimentar implementada"
#
# Reworks the meal-plan sheet:
#  - title cell text changes and its merge range shrinks by one column
#  - the per-100g macro table header/labels are renumbered (peso/G0 columns
#    removed, daily-target row values become real numbers) and the side
#    "totals" block switches from static numbers / ratio formulas to
#    SUM() formulas pulling from the food rows, plus a new "whille / <"
#    comparison row underneath
#  - workbook window is minimized

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook window state -------------------------------------------------
$excel.WindowState = -4140          # xlMinimized

# --- title row (A1:F1 -> A1:E1) --------------------------------------------
$ws.Range("A1:F1").UnMerge()
$ws.Range("A1").Value = "macros diários totais"
$ws.Range("A1:E1").Merge()

# F1 keeps being "touched" alignment-wise (no visible alignment, just the
# applyAlignment flag) now that it fell out of the merged range
$ws.Range("F1").WrapText = $true
$ws.Range("F1").WrapText = $false

# --- row 2: per-100g header / food-row 1 (frango) ---------------------------
$ws.Range("A2").ClearContents()

$ws.Range("B2").Value = "kcal"
$ws.Range("C2").Value = "prot"
$ws.Range("D2").Value = "carb"
$ws.Range("E2").Value = "gord"
$ws.Range("F2").ClearContents()

$ws.Range("L2").Value = 3.6

# --- row 3: "sedentário" daily target row -----------------------------------
$ws.Range("B3").Value = 1175
$ws.Range("C3").Value = "67.5"
$ws.Range("D3").Value = "125.0"
$ws.Range("E3").Value = "40.5"
$ws.Range("F3").ClearContents()

$ws.Range("J3").Value = 2.32
$ws.Range("K3").Value = 23.51

# --- rows 4 & 5: remaining food rows, numbers instead of text ---------------
$ws.Range("J4").Value = 2.4
$ws.Range("K4").Value = 5.2
$ws.Range("L4").Value = 0.4

$ws.Range("J5").Value = 0.9
$ws.Range("K5").Value = 2.97
$ws.Range("L5").Value = 0.15

# --- N1:Q1 totals header -----------------------------------------------------
$ws.Range("N1").Value = "kcal"
$ws.Range("N1").Font.Bold = $ws.Range("N1").Font.Bold   # no-op, keep style untouched
$ws.Range("O1").Value = "prot"
$ws.Range("P1").Value = "carb"
$ws.Range("Q1").Value = "gord"
$ws.Range("R1").ClearContents()
$ws.Range("S1").ClearContents()

# N1 loses the bold/centered style it inherited from the old "Kcal" header
$ws.Range("N1").Style = "Normal"

# --- N2:Q2 totals row: SUM() formulas instead of the old literal/ratio ones -
$ws.Range("N2").Formula = "=SUM(I2:I6)"
$ws.Range("O2").Formula = "=SUM(J2:J6)"
$ws.Range("P2").Formula = "=SUM(K2:K6)"
$ws.Range("Q2").Formula = "=SUM(L2:L6)"
$ws.Range("R2").ClearContents()
$ws.Range("S2").ClearContents()

# --- new row 9: "whille" / "<" comparison -----------------------------------
$ws.Range("N9").Value = "whille"
$ws.Range("O9").Formula = "=N2"
$ws.Range("P9").Value = "<"
$ws.Range("Q9").Formula = "=B3"
